$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (re-generated test data) ---
$ws.Range("A2").Value = "CBbjo678"
$ws.Range("B2").Value = 23110207
$ws.Range("C2").Value = "zwunmrs87"
$ws.Range("D2").Value = "s!36a%ZW"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "HhEBRHnF"
$ws.Range("G2").Value = "poJU"
$ws.Range("H2").Value = "Candidate"

# --- Update existing row 3 (re-generated test data) ---
$ws.Range("A3").Value = "GGdEM576"
$ws.Range("B3").Value = 23110206
$ws.Range("C3").Value = "zaljudy79"
$ws.Range("D3").Value = "VM!vw&79"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "ePKLMzjd"
$ws.Range("G3").Value = "Ivjd"
$ws.Range("H3").Value = "Candidate"

# --- Add a new row 4 with another generated candidate record ---
# Give the new row the same bordered look as the rows above it before
# filling in the values (mirrors the formatting already used by rows 2-3).
for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item(4, $col).Borders.LineStyle = 1
}

$ws.Range("A4").Value = "hHnMu113"
$ws.Range("B4").Value = 23110205
$ws.Range("C4").Value = "betqcob60"
$ws.Range("D4").Value = "TgH!8w5$"
$ws.Range("E4").Value = "MR"
$ws.Range("F4").Value = "NOGXBHgS"
$ws.Range("G4").Value = "oZoA"
$ws.Range("H4").Value = "Candidate"

# --- Keep the sheet's selection / used range in sync with the new row ---
$ws.Range("A1:H4").Select() | Out-Null
